# Updated computer systems spreadsheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# -----------------------------------------------------------------
# Phase 1: Insert the extra rows needed (bottom-up so row numbers
# used below always refer to the CURRENT state of the sheet at the
# time of the call).
# -----------------------------------------------------------------

# a) three new rows after the old "ESC" row (row 8) -> rows 9,10,11
$ws.Range("A9:A11").EntireRow.Insert()

# b) one new row after the old "Arduino" row (row 7) -> row 8
$ws.Range("A8:A8").EntireRow.Insert()

# c) four new rows after the "Raspberry Pi" row (row 6) -> rows 7,8,9,10
$ws.Range("A7:A10").EntireRow.Insert()

# At this point the sheet has 16 rows with this layout (A/B only):
#  1  header
#  2  Gamepad/Controller
#  3  Surface Operator Computer (start of A3:A5 merge)
#  4  (blank A)
#  5  (blank A)
#  6  Raspberry Pi
#  7  (blank A/B, new)
#  8  (blank A/B, new)
#  9  (blank A/B, new)
# 10  (blank A/B, new)
# 11  Arduino (was row 7)
# 12  (blank A/B, new)
# 13  ESC (was row 8)
# 14  (blank A/B, new)
# 15  (blank A/B, new)
# 16  (blank A/B, new)

# -----------------------------------------------------------------
# Phase 2: Fix up formatting of newly inserted rows (inserted rows
# lose their border), then add the new column C everywhere using
# PasteSpecial so every row's formatting (fill/border/alignment/
# font) is consistent with its neighbours.
# -----------------------------------------------------------------

$ws.Range("A6:B6").Copy() | Out-Null
$ws.Range("A7:B10").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("A11:B11").Copy() | Out-Null
$ws.Range("A12:B12").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("A13:B13").Copy() | Out-Null
$ws.Range("A14:B16").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# Paste the per-row formats across into column C as well (copy from
# column B only -- a single source column -- so the paste does not
# spill into column D).
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B3").Copy() | Out-Null
$ws.Range("C3:C5").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B6").Copy() | Out-Null
$ws.Range("C6:C10").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B11").Copy() | Out-Null
$ws.Range("C11:C12").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B13").Copy() | Out-Null
$ws.Range("C13:C16").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

Write-Host "Phase 1+2 done"

# -----------------------------------------------------------------
# Phase 3: wipe all cell text (formatting/merges/row heights are
# untouched) so that when we write the new text below, the shared
# string table is rebuilt cleanly in first-use order.
# -----------------------------------------------------------------

$ws.Cells.ClearContents()

# -----------------------------------------------------------------
# Phase 4: set the cell values (text), row by row, left to right.
# -----------------------------------------------------------------

$ws.Range("A1").Value2 = "System"
$ws.Range("B1").Value2 = "Tasks"
$ws.Range("C1").Value2 = "Language"

$ws.Range("A2").Value2 = "Gamepad/Controller"
$ws.Range("B2").Value2 = "Send control input to surface computer"
$ws.Range("C2").Value2 = "N/A"

$ws.Range("A3").Value2 = "Surface Operator Computer"
$ws.Range("B3").Value2 = "Relay control input to onboard computer (RasPi)"
$ws.Range("C3").Value2 = "Python"

$ws.Range("B4").Value2 = "Send commands to onboard computer"

$ws.Range("B5").Value2 = "Receive & display telemetry (camera feed, sensor data, etc.) from onboard computer"

$ws.Range("A6").Value2 = "Raspberry Pi"
$ws.Range("B6").Value2 = "Receive commands from surface, send commands to Arduino  "
$ws.Range("C6").Value2 = "Python"

$ws.Range("B7").Value2 = "Receive controller input from surface, interpret controller input, send commands to Arduino"

$ws.Range("B8").Value2 = "Receive sensor data from Arduino, send sensor data to surface"

$ws.Range("B9").Value2 = "Process sensor data, use sensor data to control vehicle (adjust camera angle, stabilize attitude, etc.)"

$ws.Range("B10").Value2 = "Receive camera feed from camera, stream camera feed to surface"

$ws.Range("A11").Value2 = "Arduino"
$ws.Range("B11").Value2 = "Read sensor data, send sensor data to onboard computer"
$ws.Range("C11").Value2 = "C++"

$ws.Range("B12").Value2 = "Receive commands from onboard computer, send control signals to ESCs, servos, etc."

$ws.Range("A13").Value2 = "ESC (Electronic Speed Controller)"
$ws.Range("B13").Value2 = "Relay motor control signals from Arduino to motor"
$ws.Range("C13").Value2 = "N/A"

$ws.Range("A14").Value2 = "6-axis sensor (Gyroscope & Accelerometer)"
$ws.Range("B14").Value2 = "Send acceleration & rotation data to Arduino"
$ws.Range("C14").Value2 = "N/A"

$ws.Range("A15").Value2 = "Thermistor (Temperature sensor)"
$ws.Range("B15").Value2 = "Send temperature data to Arduino"
$ws.Range("C15").Value2 = "N/A"

$ws.Range("A16").Value2 = "Camera"
$ws.Range("B16").Value2 = "Send video feed to onboard computer"
$ws.Range("C16").Value2 = "N/A"

Write-Host "Phase 4 done"
